$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw data values (C and D columns); the Recovery rate formulas
# in column E (D/C) will recalculate automatically.
$ws.Range("C2").Value = 1784323174
$ws.Range("D2").Value = 464050648

$ws.Range("C3").Value = 7310908806
$ws.Range("D3").Value = 1733080585

$ws.Range("C4").Value = 4253229544
$ws.Range("D4").Value = 920744754

$ws.Range("C5").Value = 7196293897
$ws.Range("D5").Value = 1504268649

$ws.Range("C6").Value = 6374400927
$ws.Range("D6").Value = 1318101869

# Update the active selection to match the new cursor position saved in the file.
$ws.Range("H13").Select()
